$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 808-865 per the diff (shift-by-2 re-shuffle of Coliflor price records) ---
# Row 808
$ws.Range("D808").Value = 45013

# Row 809
$ws.Range("D809").Value = 45013
$ws.Range("J809").Value = 790

# Row 810
$ws.Range("D810").Value = 44610
$ws.Range("J810").Value = 1600
$ws.Range("K810").Value = 1100
$ws.Range("L810").Value = 1200
$ws.Range("M810").Value = 1150
$ws.Range("P810").Value = 1150

# Row 811
$ws.Range("D811").Value = 44610
$ws.Range("J811").Value = 610
$ws.Range("K811").Value = 900
$ws.Range("L811").Value = 900
$ws.Range("M811").Value = 900
$ws.Range("P811").Value = 900

# Row 812
$ws.Range("D812").Value = 44939
$ws.Range("J812").Value = 3400

# Row 813
$ws.Range("D813").Value = 44939
$ws.Range("J813").Value = 1600

# Row 814
$ws.Range("D814").Value = 44673
$ws.Range("J814").Value = 1600
$ws.Range("K814").Value = 900
$ws.Range("L814").Value = 1000
$ws.Range("M814").Value = 950
$ws.Range("P814").Value = 950

# Row 815
$ws.Range("D815").Value = 44673
$ws.Range("J815").Value = 790
$ws.Range("K815").Value = 700
$ws.Range("M815").Value = 700
$ws.Range("P815").Value = 700

# Row 816
$ws.Range("D816").Value = 44568
$ws.Range("J816").Value = 3400
$ws.Range("K816").Value = 750
$ws.Range("L816").Value = 800
$ws.Range("M816").Value = 775
$ws.Range("O816").Value = "Región Metropolitana"
$ws.Range("P816").Value = 775

# Row 817
$ws.Range("D817").Value = 44568
$ws.Range("I817").Value = "Segunda"
$ws.Range("J817").Value = 1060
$ws.Range("K817").Value = 650
$ws.Range("L817").Value = 700
$ws.Range("M817").Value = 675
$ws.Range("O817").Value = "Región Metropolitana"
$ws.Range("P817").Value = 675

# Row 818
$ws.Range("D818").Value = 44874
$ws.Range("J818").Value = 2200
$ws.Range("K818").Value = 450
$ws.Range("L818").Value = 500
$ws.Range("M818").Value = 477
$ws.Range("O818").Value = "Provincia de Melipilla"
$ws.Range("P818").Value = 477

# Row 819
$ws.Range("D819").Value = 44306
$ws.Range("J819").Value = 1600
$ws.Range("K819").Value = 700
$ws.Range("L819").Value = 800
$ws.Range("M819").Value = 750
$ws.Range("O819").Value = "Provincia de Chacabuco"
$ws.Range("P819").Value = 750

# Row 820
$ws.Range("D820").Value = 44620
$ws.Range("J820").Value = 970
$ws.Range("K820").Value = 1100
$ws.Range("L820").Value = 1200
$ws.Range("M820").Value = 1149
$ws.Range("O820").Value = "Región Metropolitana"
$ws.Range("P820").Value = 1149

# Row 821
$ws.Range("I821").Value = "Primera"
$ws.Range("J821").Value = 4300
$ws.Range("K821").Value = 600
$ws.Range("L821").Value = 700
$ws.Range("M821").Value = 650
$ws.Range("P821").Value = 650

# Row 822
$ws.Range("I822").Value = "Primera"
$ws.Range("J822").Value = 5200
$ws.Range("K822").Value = 600
$ws.Range("L822").Value = 700
$ws.Range("M822").Value = 650
$ws.Range("P822").Value = 650

# Row 823
$ws.Range("D823").Value = 44413
$ws.Range("I823").Value = "Segunda"
$ws.Range("J823").Value = 2500
$ws.Range("K823").Value = 500
$ws.Range("L823").Value = 500
$ws.Range("M823").Value = 500
$ws.Range("P823").Value = 500

# Row 824
$ws.Range("D824").Value = 44413
$ws.Range("I824").Value = "Segunda"
$ws.Range("K824").Value = 500
$ws.Range("L824").Value = 500
$ws.Range("M824").Value = 500
$ws.Range("O824").Value = "Región de O'Higgins"
$ws.Range("P824").Value = 500

# Row 825
$ws.Range("D825").Value = 44257
$ws.Range("I825").Value = "Primera"
$ws.Range("J825").Value = 700
$ws.Range("K825").Value = 1200
$ws.Range("L825").Value = 1200
$ws.Range("M825").Value = 1200
$ws.Range("P825").Value = 1200

# Row 826
$ws.Range("D826").Value = 44426
$ws.Range("J826").Value = 3400
$ws.Range("K826").Value = 600
$ws.Range("L826").Value = 700
$ws.Range("M826").Value = 650
$ws.Range("P826").Value = 650

# Row 827
$ws.Range("D827").Value = 44426
$ws.Range("J827").Value = 970
$ws.Range("K827").Value = 400
$ws.Range("L827").Value = 400
$ws.Range("M827").Value = 400
$ws.Range("P827").Value = 400

# Row 828
$ws.Range("D828").Value = 45008
$ws.Range("J828").Value = 1600
$ws.Range("K828").Value = 800
$ws.Range("L828").Value = 900
$ws.Range("M828").Value = 850
$ws.Range("P828").Value = 850

# Row 829
$ws.Range("D829").Value = 45008
$ws.Range("I829").Value = "Segunda"
$ws.Range("J829").Value = 790
$ws.Range("K829").Value = 700
$ws.Range("L829").Value = 700
$ws.Range("M829").Value = 700
$ws.Range("O829").Value = "Región Metropolitana"
$ws.Range("P829").Value = 700

# Row 830
$ws.Range("D830").Value = 44342
$ws.Range("J830").Value = 2400
$ws.Range("K830").Value = 500
$ws.Range("L830").Value = 600
$ws.Range("M830").Value = 550
$ws.Range("P830").Value = 550

# Row 831
$ws.Range("D831").Value = 44342
$ws.Range("I831").Value = "Primera"
$ws.Range("J831").Value = 2800
$ws.Range("K831").Value = 500
$ws.Range("L831").Value = 600
$ws.Range("M831").Value = 550
$ws.Range("O831").Value = "Región de O'Higgins"
$ws.Range("P831").Value = 550

# Row 832
$ws.Range("D832").Value = 44567
$ws.Range("J832").Value = 3400
$ws.Range("K832").Value = 750
$ws.Range("L832").Value = 800
$ws.Range("M832").Value = 775
$ws.Range("P832").Value = 775

# Row 833
$ws.Range("D833").Value = 44567
$ws.Range("J833").Value = 1060
$ws.Range("K833").Value = 650
$ws.Range("L833").Value = 700
$ws.Range("M833").Value = 675
$ws.Range("P833").Value = 675

# Row 834
$ws.Range("D834").Value = 44960
$ws.Range("J834").Value = 2500
$ws.Range("K834").Value = 800
$ws.Range("L834").Value = 900
$ws.Range("M834").Value = 850
$ws.Range("O834").Value = "Región Metropolitana"
$ws.Range("P834").Value = 850

# Row 835
$ws.Range("D835").Value = 44960
$ws.Range("J835").Value = 970
$ws.Range("K835").Value = 600
$ws.Range("L835").Value = 600
$ws.Range("M835").Value = 600
$ws.Range("O835").Value = "Región Metropolitana"
$ws.Range("P835").Value = 600

# Row 836
$ws.Range("D836").Value = 44364
$ws.Range("J836").Value = 3800
$ws.Range("K836").Value = 500
$ws.Range("L836").Value = 600
$ws.Range("M836").Value = 550
$ws.Range("P836").Value = 550

# Row 837
$ws.Range("D837").Value = 44364
$ws.Range("I837").Value = "Segunda"
$ws.Range("J837").Value = 2000
$ws.Range("K837").Value = 400
$ws.Range("L837").Value = 400
$ws.Range("M837").Value = 400
$ws.Range("O837").Value = "Región de O'Higgins"
$ws.Range("P837").Value = 400

# Row 838
$ws.Range("D838").Value = 44215
$ws.Range("I838").Value = "Primera"
$ws.Range("J838").Value = 2000
$ws.Range("K838").Value = 900
$ws.Range("L838").Value = 1000
$ws.Range("M838").Value = 950
$ws.Range("O838").Value = "Región de O'Higgins"
$ws.Range("P838").Value = 950

# Row 839
$ws.Range("D839").Value = 44677
$ws.Range("K839").Value = 700
$ws.Range("L839").Value = 800
$ws.Range("M839").Value = 750
$ws.Range("P839").Value = 750

# Row 840
$ws.Range("D840").Value = 44677
$ws.Range("J840").Value = 970
$ws.Range("K840").Value = 600
$ws.Range("L840").Value = 600
$ws.Range("M840").Value = 600
$ws.Range("P840").Value = 600

# Row 841
$ws.Range("D841").Value = 44747
$ws.Range("K841").Value = 1100
$ws.Range("L841").Value = 1200
$ws.Range("M841").Value = 1150
$ws.Range("P841").Value = 1150

# Row 842
$ws.Range("D842").Value = 44747
$ws.Range("K842").Value = 900
$ws.Range("L842").Value = 900
$ws.Range("M842").Value = 900
$ws.Range("P842").Value = 900

# Row 843
$ws.Range("D843").Value = 44771
$ws.Range("J843").Value = 3400
$ws.Range("K843").Value = 900
$ws.Range("L843").Value = 1000
$ws.Range("M843").Value = 950
$ws.Range("P843").Value = 950

# Row 844
$ws.Range("D844").Value = 44771
$ws.Range("I844").Value = "Segunda"
$ws.Range("J844").Value = 1600
$ws.Range("K844").Value = 800
$ws.Range("M844").Value = 800
$ws.Range("P844").Value = 800

# Row 845
$ws.Range("D845").Value = 44245
$ws.Range("I845").Value = "Primera"
$ws.Range("K845").Value = 1200
$ws.Range("L845").Value = 1200
$ws.Range("M845").Value = 1200
$ws.Range("P845").Value = 1200

# Row 846
$ws.Range("D846").Value = 44498
$ws.Range("J846").Value = 2800
$ws.Range("K846").Value = 600
$ws.Range("M846").Value = 700
$ws.Range("P846").Value = 700

# Row 847
$ws.Range("D847").Value = 44498
$ws.Range("J847").Value = 1600
$ws.Range("K847").Value = 500
$ws.Range("L847").Value = 500
$ws.Range("M847").Value = 500
$ws.Range("P847").Value = 500

# Row 848
$ws.Range("D848").Value = 44536
$ws.Range("J848").Value = 2500
$ws.Range("K848").Value = 700
$ws.Range("L848").Value = 800
$ws.Range("M848").Value = 750
$ws.Range("P848").Value = 750

# Row 849
$ws.Range("D849").Value = 44536
$ws.Range("K849").Value = 600
$ws.Range("L849").Value = 600
$ws.Range("M849").Value = 600
$ws.Range("P849").Value = 600

# Row 850
$ws.Range("D850").Value = 44972

# Row 851
$ws.Range("D851").Value = 44972
$ws.Range("I851").Value = "Segunda"
$ws.Range("J851").Value = 790
$ws.Range("K851").Value = 900
$ws.Range("L851").Value = 900
$ws.Range("M851").Value = 900
$ws.Range("P851").Value = 900

# Row 852
$ws.Range("D852").Value = 44671
$ws.Range("I852").Value = "Primera"
$ws.Range("J852").Value = 1600
$ws.Range("K852").Value = 1000
$ws.Range("L852").Value = 1100
$ws.Range("M852").Value = 1050
$ws.Range("P852").Value = 1050

# Row 853
$ws.Range("D853").Value = 44988
$ws.Range("J853").Value = 1600
$ws.Range("K853").Value = 1000
$ws.Range("L853").Value = 1100
$ws.Range("M853").Value = 1050
$ws.Range("P853").Value = 1050

# Row 854
$ws.Range("D854").Value = 44988
$ws.Range("J854").Value = 880
$ws.Range("K854").Value = 800
$ws.Range("L854").Value = 800
$ws.Range("M854").Value = 800
$ws.Range("O854").Value = "Región Metropolitana"
$ws.Range("P854").Value = 800

# Row 855
$ws.Range("D855").Value = 44608
$ws.Range("J855").Value = 3400
$ws.Range("K855").Value = 900
$ws.Range("L855").Value = 1000
$ws.Range("M855").Value = 950
$ws.Range("P855").Value = 950

# Row 856
$ws.Range("D856").Value = 44608
$ws.Range("I856").Value = "Segunda"
$ws.Range("J856").Value = 1060
$ws.Range("K856").Value = 700
$ws.Range("M856").Value = 700
$ws.Range("O856").Value = "Provincia de Chacabuco"
$ws.Range("P856").Value = 700

# Row 857
$ws.Range("I857").Value = "Primera"
$ws.Range("J857").Value = 4300
$ws.Range("K857").Value = 600
$ws.Range("L857").Value = 700
$ws.Range("M857").Value = 650
$ws.Range("P857").Value = 650

# Row 858
$ws.Range("I858").Value = "Primera"
$ws.Range("J858").Value = 5200
$ws.Range("K858").Value = 600
$ws.Range("L858").Value = 700
$ws.Range("M858").Value = 650
$ws.Range("P858").Value = 650

# Row 859
$ws.Range("D859").Value = 44414
$ws.Range("I859").Value = "Segunda"
$ws.Range("J859").Value = 1600
$ws.Range("K859").Value = 500
$ws.Range("L859").Value = 500
$ws.Range("M859").Value = 500
$ws.Range("P859").Value = 500

# Row 860
$ws.Range("D860").Value = 44414
$ws.Range("J860").Value = 2500
$ws.Range("K860").Value = 500
$ws.Range("L860").Value = 500
$ws.Range("M860").Value = 500
$ws.Range("O860").Value = "Región de O'Higgins"
$ws.Range("P860").Value = 500

# Row 861
$ws.Range("D861").Value = 44925
$ws.Range("J861").Value = 3400
$ws.Range("K861").Value = 900
$ws.Range("L861").Value = 1000
$ws.Range("M861").Value = 950
$ws.Range("O861").Value = "Región Metropolitana"
$ws.Range("P861").Value = 950

# Row 862
$ws.Range("D862").Value = 44925
$ws.Range("I862").Value = "Segunda"
$ws.Range("K862").Value = 700
$ws.Range("L862").Value = 700
$ws.Range("M862").Value = 700
$ws.Range("P862").Value = 700

# Row 863
$ws.Range("D863").Value = 44189
$ws.Range("I863").Value = "Primera"
$ws.Range("J863").Value = 2500
$ws.Range("K863").Value = 800
$ws.Range("M863").Value = 840
$ws.Range("O863").Value = "Región de O'Higgins"
$ws.Range("P863").Value = 840

# Row 864
$ws.Range("D864").Value = 44601
$ws.Range("J864").Value = 1600
$ws.Range("K864").Value = 1000
$ws.Range("L864").Value = 1200

# Row 865
$ws.Range("D865").Value = 44601

# --- Append two brand-new rows (866, 867) at the end of the data block ---
# Row 866
$ws.Range("A866").Value = 9
$ws.Range("B866").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C866").Value = "Metropolitana"
$ws.Range("D866").Value = 44839
$ws.Range("E866").Value = 13
$ws.Range("F866").Value = 100112008
$ws.Range("G866").Value = "Coliflor"
$ws.Range("H866").Value = "Sin especificar"
$ws.Range("I866").Value = "Primera"
$ws.Range("J866").Value = 1060
$ws.Range("K866").Value = 1100
$ws.Range("L866").Value = 1100
$ws.Range("M866").Value = 1100
$ws.Range("N866").Value = "`$/unidad"
$ws.Range("O866").Value = "Región Metropolitana"
$ws.Range("P866").Value = 1100
$ws.Range("Q866").Value = 1
$ws.Range("R866").Value = "Hortaliza"
$ws.Range("D866").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 867
$ws.Range("A867").Value = 9
$ws.Range("B867").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C867").Value = "Metropolitana"
$ws.Range("D867").Value = 44839
$ws.Range("E867").Value = 13
$ws.Range("F867").Value = 100112008
$ws.Range("G867").Value = "Coliflor"
$ws.Range("H867").Value = "Sin especificar"
$ws.Range("I867").Value = "Segunda"
$ws.Range("J867").Value = 610
$ws.Range("K867").Value = 900
$ws.Range("L867").Value = 900
$ws.Range("M867").Value = 900
$ws.Range("N867").Value = "`$/unidad"
$ws.Range("O867").Value = "Región Metropolitana"
$ws.Range("P867").Value = 900
$ws.Range("Q867").Value = 1
$ws.Range("R867").Value = "Hortaliza"
$ws.Range("D867").NumberFormat = "YYYY-MM-DD HH:MM:SS"
